$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column headers: J = "MST" (algorithm name), K = "%Over" (relative gap vs optimal)
$ws.Range("J1").Value = "MST"
$ws.Range("K1").Value = "%Over"
$ws.Range("K1").NumberFormat = $ws.Range("I1").NumberFormat

# MST 2-approximation results (column J) and computed %Over (column K = J/B - 1)
$mstValues = @(2380448, 10402, 1150963, 65712, 301216, 134748, 2027107, 1646249, 838282, 1134989, 1675105, 68090, 178249)

for ($i = 0; $i -lt $mstValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 10).Value = $mstValues[$i]
    $ws.Cells.Item($row, 11).Formula = "=J$row/B$row-1"
    $ws.Cells.Item($row, 11).NumberFormat = $ws.Range("I1").NumberFormat
}

[void]$ws.Range("J15").Select()
